$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.206503629684448
$ws.Range("B1").Value = 2.575319766998291
$ws.Range("C1").Value = 9.136944770812988
$ws.Range("D1").Value = 2.033324718475342
$ws.Range("E1").Value = 1.166728258132935
